# RuleToJson.xlsx - "Add files via upload" edit
#
# 1. Rule id in row 2 changes from 32 to 1 (A2: 32 -> 1).
# 2. Both rule-building formulas (R2 and R3) are updated so the "first
#    rule" check also accepts a text "1" in column A, not just the
#    number 1: IF(A#=1,...) -> IF(OR(A#=1,A#="1"),...)
# 3. The active selection on Sheet1 moves from J16 to I7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. A2: 32 -> 1 -------------------------------------------------
$ws.Range("A2").Value = 1

# --- 2. Patch the IF(...) guard in R2 and R3 -------------------------
# The formulas are long, so read back the existing text and only touch
# the one clause that changed, rather than retyping the whole formula.
# Embedded line breaks are stripped first: re-assigning a `.Formula`
# that still contains literal CR/LF makes this host auto-fit the row,
# which the real edit never did (no row-height change is present).

$r2Formula = $ws.Range("R2").Formula
$r2Formula = $r2Formula.Replace([char]13, "").Replace([char]10, "")
$r2Formula = $r2Formula.Replace("IF(A2=1,", "IF(OR(A2=1,A2=`"1`"),")
$ws.Range("R2").Formula = $r2Formula

$r3Formula = $ws.Range("R3").Formula
$r3Formula = $r3Formula.Replace([char]13, "").Replace([char]10, "")
$r3Formula = $r3Formula.Replace("IF(A3=1,", "IF(OR(A3=1,A3=`"1`"),")
$ws.Range("R3").Formula = $r3Formula

# --- 3. Selection: J16 -> I7 -----------------------------------------
$ws.Range("I7").Select()
